$d = $word.ActiveDocument

$pairs = @(
    @("2024-03-18 Monday", "2024-03-19 Tuesday"),
    @("60×66=", "14×20="),
    @("73×81=", "28×51="),
    @("99×83=", "15×99="),
    @("19×28=", "47×48="),
    @("61×53=", "91×77="),
    @("72×38=", "89×34="),
    @("86×26=", "38×75="),
    @("40×22=", "92×34="),
    @("48×16=", "30×62="),
    @("79×73=", "20×36="),
    @("98×38=", "39×57="),
    @("81×57=", "20×64="),
    @("43×99=", "56×78="),
    @("42×28=", "57×16="),
    @("44×82=", "68×85="),
    @("71×11=", "38×53="),
    @("79×33=", "93×24="),
    @("62×39=", "93×89="),
    @("73×80=", "25×61="),
    @("39×90=", "57×90="),
    @("51×50=", "86×45="),
    @("89×79=", "48×58="),
    @("39×91=", "87×28="),
    @("41×96=", "20×46="),
    @("13×49=", "14×12=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
